$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update 想去人数 (column F) for several rows
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 521
$ws1.Range("F6").Value = 6946
$ws1.Range("F7").Value = 195
$ws1.Range("F8").Value = 157
$ws1.Range("F9").Value = 1048
$ws1.Range("F10").Value = 410
$ws1.Range("F11").Value = 143
$ws1.Range("F13").Value = 591

# Sheet "演出" (sheet2): update 想去人数 (column F)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 20

# Sheet "全部类型" (sheet4): update 想去人数 (column F) for several rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 521
$ws4.Range("F6").Value = 6946
$ws4.Range("F7").Value = 195
$ws4.Range("F8").Value = 157
$ws4.Range("F9").Value = 1048
$ws4.Range("F10").Value = 410
$ws4.Range("F11").Value = 143
$ws4.Range("F13").Value = 591
$ws4.Range("F14").Value = 20
